$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "22.407.64"
$ws.Range("E2").Value = "  -0.31%  "

# Row 3
$ws.Range("D3").Value = "1.567.27"
$ws.Range("E3").Value = "  -0.18%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.000"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.13%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "286.87"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.24%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3746"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +2.66%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3275"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.64%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.52"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -5.38%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.152"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.19%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07423"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.04%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.03%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.51"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.69%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.846"
$ws.Range("D14").ClearFormats()

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.835"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.11%  "

# Row 16
$ws.Range("D16").Value = "1.563.83"
$ws.Range("E16").Value = "  -0.33%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001100"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.51%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06701"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.46%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "85.86"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.32%  "

# Row 20
$ws.Range("E20").Value = "  -0.13%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.355"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.35%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.28"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.84%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.70"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.79%  "

# Row 24
$ws.Range("D24").Value = "22.417.87"
$ws.Range("E24").Value = "  -0.24%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.308"
$ws.Range("D25").ClearFormats()

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.556"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.96%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "151.00"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.27%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.41"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.61%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.915"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.99%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.44"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.61%  "

# Row 31
$ws.Range("D31").Value = "1.741.15"
$ws.Range("E31").Value = "  -0.19%  "

# Row 32
$ws.Range("E32").Value = "  +2.28%  "

# Row 33
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.930"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.32%  "

# Row 34
$ws.Range("B34").Value = "WEMIXTOKEN"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.941"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.58%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.654"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.21%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08228"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.37%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02391"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.37%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.298"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.05%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06323"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.94%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2190"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.43%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.264"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.72%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.12"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.67%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6107"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.78%  "

# Row 44
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.79"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.71%  "

# Row 45
$ws.Range("B45").Value = "PancakeSwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.750"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.06%  "

# Row 46
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5921"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.88%  "

# Row 47
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.010"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.73%  "

# Row 48
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "123.75"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.38%  "

# Row 49
$ws.Range("B49").Value = "EOS"
$ws.Range("C49").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.179"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.88%  "

# Row 50
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07151"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.79%  "

# Row 51
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "75.82"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.67%  "
